$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '40.483.55'
$ws.Range('E2').Value = '  +1.00%  '

# Row 3
$ws.Range('D3').Value = '2.227.10'
$ws.Range('E3').Value = '  -0.16%  '

# Row 4
$ws.Range('E4').Value = '  -0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '89.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.81%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.517'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.83%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.18%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.475'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.54%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.06%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '31.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.64%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0788'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.94%  '

# Row 13
$ws.Range('E13').Value = '  +2.58%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.33%  '

# Row 15
$ws.Range('D15').Value = '2.569.62'
$ws.Range('E15').Value = '  -0.13%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.94'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.94%  '

# Row 17
$ws.Range('D17').Value = '2.242.27'
$ws.Range('E17').Value = '  +0.60%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.742'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.75%  '

# Row 19
$ws.Range('D19').Value = '40.423.78'
$ws.Range('E19').Value = '  +1.03%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0893'
$ws.Range('E20').Value = '  +0.66%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.36%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '237.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.30%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.03%  '

# Row 26
$ws.Range('E26').Value = '  -0.08%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.51%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.45%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.87%  '

# Row 30
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.04%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '156.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.49%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.54%  '

# Row 33
$ws.Range('E33').Value = '  +0.07%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.00%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0723'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.78%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.32%  '

# Row 37
$ws.Range('E37').Value = '  -0.55%  '

# Row 38
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.103'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.15%  '

# Row 39
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.114'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.77%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.32%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.78'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.31%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.62%  '

# Row 43
$ws.Range('D43').Value = '2.077.81'
$ws.Range('E43').Value = '  -2.13%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.54%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0273'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.06%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.12'
$ws.Range('D46').Style = 'Normal'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.32%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.85'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -13.13%  '

# Row 49
$ws.Range('D49').Value = '2.438.17'
$ws.Range('E49').Value = '  +0.15%  '

# Row 50
$ws.Range('E50').Value = '  +2.49%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.38%  '
